$wb = $excel.ActiveWorkbook

# stage through unique temp names first to avoid name collisions
$wb.Worksheets.Item(1).Name = "__tmp_sheet_1__"
$wb.Worksheets.Item(2).Name = "__tmp_sheet_2__"
$wb.Worksheets.Item(3).Name = "__tmp_sheet_3__"
$wb.Worksheets.Item(4).Name = "__tmp_sheet_4__"
$wb.Worksheets.Item(5).Name = "__tmp_sheet_5__"
$wb.Worksheets.Item(6).Name = "__tmp_sheet_6__"
$wb.Worksheets.Item(7).Name = "__tmp_sheet_7__"
$wb.Worksheets.Item(8).Name = "__tmp_sheet_8__"
$wb.Worksheets.Item(9).Name = "__tmp_sheet_9__"

# ---- sheet 1 ----
$ws = $wb.Worksheets.Item(1)
$ws.Name = "summ3"
$ws.Rows.Item(5).Delete()

$data = @(
  @("Intercept", 1.813174806798406, 0.005104391834953413),
  @("HHType_simp[T.MultiAdult_Kids]", 0.08344541531016812, 0.8025015870775489),
  @("HHType_simp[T.Single_Female]", -0.9502765436374905, 0.000000506489104903999),
  @("HHType_simp[T.Single_Male]", -0.8626577484887158, 0.000009004023125104923),
  @("HHType_simp[T.Single_Parent]", -1.065318792408792, 0.0000003999890568030412),
  @("HHSize", 0.0238249556205674, 0.8503324805206582),
  @("IncomeDetailed_Numeric", 0.0008950411732451256, 0.0000000000000000000000000000000000000992124434180629),
  @("maxAgeHH", -0.01089200748048978, 0.00109286766608866),
  @("UniversityEducation", 0.5813962405755028, 0.0000002016185532902557),
  @("InEmployment", 0.8130919872209085, 0.00000000009127692971165613),
  @("AllRetired", 0.4190505469608328, 0.0163894393632822),
  @("UrbPopDensity", -0.00005524183628065876, 0.108945257087465),
  @("UrbBuildDensity", -0.000000001212647396431365, 0.9723585751489892),
  @("DistSubcenter", -0.03494523266907492, 0.318845553079817),
  @("DistCenter", 0.03381281514527027, 0.09156221574747328),
  @("bike_lane_share", -2.147790325649471, 0.0000009511861440370524),
  @("IntersecDensity", -0.006689850331062489, 0.03901661917909389),
  @("StreetLength", -0.002362155020529284, 0.4889582909745394),
  @("LU_UrbFab", 4.589083288612381, 0.000000307256618584674),
  @("LU_Comm", 4.404402181072343, 0.000005256256710765216),
  @("LU_Urban", -3.772350148190748, 0.000006198156556309899)
)

$r = 2
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $r = $r + 1
}

# ---- sheet 2 ----
$ws = $wb.Worksheets.Item(2)
$ws.Name = "summ6"
$ws.Rows.Item(5).Delete()

$data = @(
  @("Intercept", 1.565783875185682, 0.01709487733296084),
  @("HHType_simp[T.MultiAdult_Kids]", 0.0009685071539861091, 0.9975996425587864),
  @("HHType_simp[T.Single_Female]", -0.9757864100965854, 0.0000001902865079097114),
  @("HHType_simp[T.Single_Male]", -0.7624867616448054, 0.00007938988741677851),
  @("HHType_simp[T.Single_Parent]", -0.8798550530051668, 0.00003243322416679518),
  @("HHSize", 0.003114318123036901, 0.9800828724867547),
  @("IncomeDetailed_Numeric", 0.0009282194522482169, 0.000000000000000000000000000000000000001714328436081758),
  @("maxAgeHH", -0.01125454072926269, 0.0007732259287870918),
  @("UniversityEducation", 0.620141642401454, 0.00000003240845009079462),
  @("InEmployment", 0.6513770351228277, 0.0000002339466618474136),
  @("AllRetired", 0.3483045222368187, 0.04773342047018968),
  @("UrbPopDensity", -0.00005585317072439002, 0.1040432866146453),
  @("UrbBuildDensity", -0.000000009435391038555825, 0.7843099151190677),
  @("DistSubcenter", -0.03154410209787165, 0.3750189970890007),
  @("DistCenter", 0.04598526709392593, 0.02546250319552558),
  @("bike_lane_share", -2.009157662310828, 0.000004661879417062339),
  @("IntersecDensity", -0.004480863091115468, 0.1670808570942043),
  @("StreetLength", -0.0004666313292058907, 0.895395780735512),
  @("LU_UrbFab", 3.718672140990676, 0.00002757332034360551),
  @("LU_Comm", 3.271962843060712, 0.0006355553033106536),
  @("LU_Urban", -3.104717277861679, 0.0001945685466861046)
)

$r = 2
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $r = $r + 1
}

# ---- sheet 3 ----
$ws = $wb.Worksheets.Item(3)
$ws.Name = "summ1"
$ws.Rows.Item(5).Delete()

$data = @(
  @("Intercept", 1.387184252976138, 0.03258097237474563),
  @("HHType_simp[T.MultiAdult_Kids]", 0.1607581396822209, 0.6159726114168804),
  @("HHType_simp[T.Single_Female]", -0.8427262671314113, 0.000005491190144914988),
  @("HHType_simp[T.Single_Male]", -0.6673735199257992, 0.0004749216525973628),
  @("HHType_simp[T.Single_Parent]", -0.884312460403263, 0.00002481591527522136),
  @("HHSize", -0.008456633578425934, 0.9448644563737841),
  @("IncomeDetailed_Numeric", 0.0008892394162480711, 0.00000000000000000000000000000000000003724209076716121),
  @("maxAgeHH", -0.005829403254286632, 0.07883700599114134),
  @("UniversityEducation", 0.636278689365345, 0.000000008905235174613364),
  @("InEmployment", 0.7514271007430828, 0.000000001452072102594885),
  @("AllRetired", 0.2663390709605743, 0.1274552505217131),
  @("UrbPopDensity", -0.00005481690525255317, 0.1101853154093256),
  @("UrbBuildDensity", 0.00000001895915066888887, 0.5862315933517547),
  @("DistSubcenter", -0.04777048761937593, 0.1715663587231055),
  @("DistCenter", 0.04368215522939035, 0.02993482942194262),
  @("bike_lane_share", -2.103788832806742, 0.000002168610500907534),
  @("IntersecDensity", -0.007177694973363112, 0.02564373983075239),
  @("StreetLength", -0.0007980770158964811, 0.8220936469795013),
  @("LU_UrbFab", 3.58787539957818, 0.00005275878164383847),
  @("LU_Comm", 2.924497216096088, 0.002215415866580119),
  @("LU_Urban", -3.002362142107893, 0.0002849606264002294)
)

$r = 2
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $r = $r + 1
}

# ---- sheet 4 ----
$ws = $wb.Worksheets.Item(4)
$ws.Name = "summ0"
$ws.Rows.Item(5).Delete()

$data = @(
  @("Intercept", 1.878342335654579, 0.004521341994143597),
  @("HHType_simp[T.MultiAdult_Kids]", 0.2632864117744441, 0.4283461487289144),
  @("HHType_simp[T.Single_Female]", -0.8679987606092295, 0.000005660059730310935),
  @("HHType_simp[T.Single_Male]", -0.7128669253143962, 0.0002996842232792962),
  @("HHType_simp[T.Single_Parent]", -0.9345474320382953, 0.00000663470760743353),
  @("HHSize", -0.01714100716259191, 0.8961935354795173),
  @("IncomeDetailed_Numeric", 0.0009083787313151873, 0.000000000000000000000000000000000000002408803738617958),
  @("maxAgeHH", -0.008953923721124206, 0.007443067137626981),
  @("UniversityEducation", 0.6037796501132942, 0.00000005078761629606387),
  @("InEmployment", 0.727985915598375, 0.000000007069023220340713),
  @("AllRetired", 0.3166530700296515, 0.06996831377358399),
  @("UrbPopDensity", -0.00005523440575330044, 0.1102818358018725),
  @("UrbBuildDensity", 0.000000009201560504528055, 0.7948752287631995),
  @("DistSubcenter", -0.05444760121544806, 0.1213412326342785),
  @("DistCenter", 0.05068343958747529, 0.01316497072240227),
  @("bike_lane_share", -1.799392904708549, 0.00003916101551573561),
  @("IntersecDensity", -0.007700536097564215, 0.01709911601672692),
  @("StreetLength", -0.002683541284870736, 0.4573543193578146),
  @("LU_UrbFab", 4.318263560665627, 0.000001143964855980598),
  @("LU_Comm", 3.771575703619984, 0.00007285483907079683),
  @("LU_Urban", -3.742925260952057, 0.000006098743840846794)
)

$r = 2
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $r = $r + 1
}

# ---- sheet 5 ----
$ws = $wb.Worksheets.Item(5)
$ws.Name = "summ8"
$ws.Rows.Item(5).Delete()

$data = @(
  @("Intercept", 1.859075000983697, 0.004831717742597535),
  @("HHType_simp[T.MultiAdult_Kids]", 0.274515276361752, 0.3983630834081353),
  @("HHType_simp[T.Single_Female]", -0.9702521165713696, 0.0000001457210949291738),
  @("HHType_simp[T.Single_Male]", -0.785704110944658, 0.00003628885857531557),
  @("HHType_simp[T.Single_Parent]", -0.8636248276118563, 0.00003918964300530722),
  @("HHSize", -0.04226383751918579, 0.7286051993906528),
  @("IncomeDetailed_Numeric", 0.0009039166589027113, 0.000000000000000000000000000000000000006115576068303767),
  @("maxAgeHH", -0.01016883919041443, 0.002016574023734333),
  @("UniversityEducation", 0.6166961669575062, 0.00000002785890109882327),
  @("InEmployment", 0.6576309139508465, 0.0000001278342523974316),
  @("AllRetired", 0.3997152730880519, 0.02145153124892551),
  @("UrbPopDensity", -0.00003483224190210311, 0.3081254964245593),
  @("UrbBuildDensity", -0.00000001582476032313192, 0.6436100553533779),
  @("DistSubcenter", -0.0172698724414879, 0.6190510612699975),
  @("DistCenter", 0.0360749877408282, 0.07150462296696468),
  @("bike_lane_share", -1.913330490332726, 0.000009370127732823579),
  @("IntersecDensity", -0.007272851759850706, 0.02434802449389608),
  @("StreetLength", -0.002921276624899535, 0.4176377432252515),
  @("LU_UrbFab", 3.166785938951158, 0.0003775816335200737),
  @("LU_Comm", 3.203719128071288, 0.0008542327706604719),
  @("LU_Urban", -2.74118405098837, 0.0009667241946657101)
)

$r = 2
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $r = $r + 1
}

# ---- sheet 6 ----
$ws = $wb.Worksheets.Item(6)
$ws.Name = "summ11"
$ws.Rows.Item(5).Delete()

$data = @(
  @("Intercept", 1.941884670395376, 0.002446607422699896),
  @("HHType_simp[T.MultiAdult_Kids]", 0.2802693738529037, 0.3941170180305159),
  @("HHType_simp[T.Single_Female]", -0.9739170816976317, 0.00000014169155684124),
  @("HHType_simp[T.Single_Male]", -0.908803906350575, 0.000001713049943892752),
  @("HHType_simp[T.Single_Parent]", -0.8775041054294036, 0.00003322675518417612),
  @("HHSize", -0.07442406416355679, 0.5366935854166255),
  @("IncomeDetailed_Numeric", 0.0008616633214419922, 0.00000000000000000000000000000000001745510559736454),
  @("maxAgeHH", -0.007675383672555891, 0.02174928139896592),
  @("UniversityEducation", 0.6839365398830305, 0.0000000009135692240024245),
  @("InEmployment", 0.7330618401469788, 0.000000005418335873469133),
  @("AllRetired", 0.3454191258498391, 0.05106427288462257),
  @("UrbPopDensity", -0.00006172274738663715, 0.07633830531131879),
  @("UrbBuildDensity", 0.00000002088189257292385, 0.5539779612820075),
  @("DistSubcenter", -0.066309278696608, 0.06012876689971022),
  @("DistCenter", 0.05310982770803135, 0.009873133907030806),
  @("bike_lane_share", -1.927510834464014, 0.00001028930689489767),
  @("IntersecDensity", -0.007933497491248312, 0.0134910464296705),
  @("StreetLength", -0.002577706456781253, 0.4448538631403028),
  @("LU_UrbFab", 4.142987219316782, 0.000004302703538904852),
  @("LU_Comm", 3.357313790923919, 0.0005840985233270355),
  @("LU_Urban", -3.3716683701274, 0.00006116005069545997)
)

$r = 2
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $r = $r + 1
}

# ---- sheet 7 ----
$ws = $wb.Worksheets.Item(7)
$ws.Name = "summ5"
$ws.Rows.Item(5).Delete()

$data = @(
  @("Intercept", 1.449560191455598, 0.02591068052280267),
  @("HHType_simp[T.MultiAdult_Kids]", 0.1345889448775641, 0.6712053609646585),
  @("HHType_simp[T.Single_Female]", -0.9800776494665505, 0.00000008464886455478531),
  @("HHType_simp[T.Single_Male]", -0.7599763257463765, 0.00005228310677495268),
  @("HHType_simp[T.Single_Parent]", -1.001487686443425, 0.000001131030185794503),
  @("HHSize", -0.03184422301999832, 0.7882201937093651),
  @("IncomeDetailed_Numeric", 0.0008768758151306486, 0.0000000000000000000000000000000000001396849232938283),
  @("maxAgeHH", -0.00726908423132708, 0.02985570408353964),
  @("UniversityEducation", 0.6244341425575153, 0.00000001898487314171218),
  @("InEmployment", 0.7266180301589046, 0.000000006163198651467403),
  @("AllRetired", 0.2893584609034289, 0.100630889252375),
  @("UrbPopDensity", -0.00007877579005588995, 0.02033277215018855),
  @("UrbBuildDensity", 0.0000000006935553841620705, 0.9838209717651452),
  @("DistSubcenter", -0.0561351701139961, 0.1083291686742911),
  @("DistCenter", 0.04441872350990467, 0.02988010075767602),
  @("bike_lane_share", -1.624895955388881, 0.0002065281901686521),
  @("IntersecDensity", -0.004194098688910193, 0.1977501925446381),
  @("StreetLength", 0.0001627405101375958, 0.96438856607888),
  @("LU_UrbFab", 3.717765132853278, 0.00002872016775387366),
  @("LU_Comm", 2.910876290627638, 0.002286178287006878),
  @("LU_Urban", -3.000058109166132, 0.0002923346785128166)
)

$r = 2
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $r = $r + 1
}

# ---- sheet 8 ----
$ws = $wb.Worksheets.Item(8)
$ws.Name = "summ18"
$ws.Rows.Item(5).Delete()

$data = @(
  @("Intercept", 2.041614475544934, 0.001771384766418845),
  @("HHType_simp[T.MultiAdult_Kids]", 0.1608450351290502, 0.6220525102376362),
  @("HHType_simp[T.Single_Female]", -0.9842840709257696, 0.0000001376520590729936),
  @("HHType_simp[T.Single_Male]", -0.7273035436088386, 0.0001692944061930148),
  @("HHType_simp[T.Single_Parent]", -0.8063040423300353, 0.0001899310457905731),
  @("HHSize", -0.02496116432413896, 0.841707549748282),
  @("IncomeDetailed_Numeric", 0.0008749859244322903, 0.0000000000000000000000000000000000007847123698807711),
  @("maxAgeHH", -0.009530729102457365, 0.004440302491696252),
  @("UniversityEducation", 0.6605149844985304, 0.000000003416215034476677),
  @("InEmployment", 0.6656907225811025, 0.0000001027830913180144),
  @("AllRetired", 0.2715918603538913, 0.1215796870999206),
  @("UrbPopDensity", -0.00005953780296394347, 0.0847286327443143),
  @("UrbBuildDensity", 0.000000009844805421249736, 0.780845569328156),
  @("DistSubcenter", -0.04271853086740748, 0.2275398890322026),
  @("DistCenter", 0.05238927804852456, 0.01136871722418901),
  @("bike_lane_share", -1.94654068459514, 0.000008290390433960836),
  @("IntersecDensity", -0.006256245786118954, 0.05206849296074251),
  @("StreetLength", -0.003020915317806658, 0.3795923443862783),
  @("LU_UrbFab", 4.247250541819637, 0.000001517481011343489),
  @("LU_Comm", 3.563560022352247, 0.0001617985387513607),
  @("LU_Urban", -3.710074734893025, 0.000006578644789055394)
)

$r = 2
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $r = $r + 1
}

# ---- sheet 9 ----
$ws = $wb.Worksheets.Item(9)
$ws.Name = "summ14"
$ws.Rows.Item(5).Delete()

$data = @(
  @("Intercept", 1.922275694347778, 0.003073341324957174),
  @("HHType_simp[T.MultiAdult_Kids]", 0.1027692756898404, 0.7520093028023809),
  @("HHType_simp[T.Single_Female]", -0.9063702969807796, 0.000001394086053966768),
  @("HHType_simp[T.Single_Male]", -0.771997755106432, 0.00006529326185783096),
  @("HHType_simp[T.Single_Parent]", -0.9673987618278456, 0.000003265151253932049),
  @("HHSize", 0.00980091010240055, 0.9375142933768436),
  @("IncomeDetailed_Numeric", 0.0008739438064722794, 0.000000000000000000000000000000000003102416826620219),
  @("maxAgeHH", -0.008537045331768798, 0.01007719944230847),
  @("UniversityEducation", 0.6408213529333323, 0.000000009431128502870621),
  @("InEmployment", 0.7290039244982485, 0.000000006202356212658178),
  @("AllRetired", 0.3042350368837174, 0.08115019047276359),
  @("UrbPopDensity", -0.00004389899188780409, 0.1996310221331261),
  @("UrbBuildDensity", -0.000000007338798812120433, 0.8325871346626303),
  @("DistSubcenter", -0.02984998266353509, 0.3965627000287086),
  @("DistCenter", 0.03276617559908143, 0.1052339839045828),
  @("bike_lane_share", -2.216567721075259, 0.0000004841123117277633),
  @("IntersecDensity", -0.007193020673798084, 0.02416141351723627),
  @("StreetLength", -0.002381922689768958, 0.4866617016376846),
  @("LU_UrbFab", 3.779820598565719, 0.00002069690340512314),
  @("LU_Comm", 3.675767837207787, 0.0001216030577658463),
  @("LU_Urban", -3.374709943698205, 0.000045663837183769)
)

$r = 2
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $r = $r + 1
}

